$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the severity legend cells from columns C/D to columns E/F on rows 2-4 ---
# Row 2: C2 -> E2, D2 -> F2
$ws.Range("D2").Copy($ws.Range("F2"))
$ws.Range("D2").Clear()
$ws.Range("C2").Copy($ws.Range("E2"))
$ws.Range("C2").Clear()

# Row 3: D3 -> F3
$ws.Range("D3").Copy($ws.Range("F3"))
$ws.Range("D3").Clear()

# Row 4: D4 -> F4
$ws.Range("D4").Copy($ws.Range("F4"))
$ws.Range("D4").Clear()

# --- Shift the target/pentester labels down one slot and insert new fields ---
# A4 (Pentester Name) moves to C4
$ws.Range("A4").Copy($ws.Range("C4"))
# A3 (Target Name) moves to A4
$ws.Range("A3").Copy($ws.Range("A4"))

# New header fields
$ws.Range("A3").Value = "Target URL/IP : 192.168.1.4:3000"
$ws.Range("C3").Value = "Start Date : 2023-05-14"

# --- Merge the newly introduced label cells across A:B ---
$ws.Range("A3:B3").Merge()
$ws.Range("A4:B4").Merge()

Write-Output "edit complete"
